$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames (sharedStrings): ht_goals_h -> HTHG, ht_goals_a -> HTAG
$ws.Range("I1").Value = "HTHG"
$ws.Range("J1").Value = "HTAG"

# Row data corrections (match results / odds re-synced against upstream source)

# Row 235
$ws.Range("B235").Value = 6865915
$ws.Range("E235").Value = "FC Voluntari"
$ws.Range("F235").Value = "Universitatea Cluj"
$ws.Range("G235").Value = 0
$ws.Range("H235").Value = 0
$ws.Range("J235").Value = 0
$ws.Range("K235").Value = "D"
$ws.Range("L235").Value = 3.5
$ws.Range("M235").Value = 3.25
$ws.Range("N235").Value = 2.05
$ws.Range("O235").Value = 3.4
$ws.Range("P235").Value = 3.1
$ws.Range("Q235").Value = 2.15
$ws.Range("S235").Value = 1.975
$ws.Range("T235").Value = 1.875
$ws.Range("V235").Value = 2.05
$ws.Range("W235").Value = 1.75
$ws.Range("Y235").Value = 2.1
$ws.Range("Z235").Value = -1
$ws.Range("AA235").Value = 0.4875
$ws.Range("AB235").Value = -0.5
$ws.Range("AC235").Value = -1
$ws.Range("AD235").Value = 0.75

# Row 236
$ws.Range("B236").Value = 6861095
$ws.Range("E236").Value = "FC Botosani"
$ws.Range("F236").Value = "Farul Constanta"
$ws.Range("L236").Value = 3.75
$ws.Range("M236").Value = 3.4
$ws.Range("N236").Value = 1.909
$ws.Range("O236").Value = 3.1
$ws.Range("P236").Value = 3
$ws.Range("Q236").Value = 2.375
$ws.Range("S236").Value = 1.775
$ws.Range("T236").Value = 2.1
$ws.Range("U236").Value = 2
$ws.Range("V236").Value = 1.8
$ws.Range("W236").Value = 2.05
$ws.Range("Y236").Value = 2
$ws.Range("AA236").Value = 0.3875
$ws.Range("AD236").Value = 1.05

# Row 238
$ws.Range("B238").Value = 6870268
$ws.Range("E238").Value = "Petrolul Ploiesti"
$ws.Range("F238").Value = "ACS Sepsi"
$ws.Range("H238").Value = 2
$ws.Range("J238").Value = 1
$ws.Range("K238").Value = "A"
$ws.Range("L238").Value = 2.8
$ws.Range("M238").Value = 3
$ws.Range("N238").Value = 2.55
$ws.Range("O238").Value = 3
$ws.Range("P238").Value = 3.2
$ws.Range("Q238").Value = 2.3
$ws.Range("R238").Value = 0.25
$ws.Range("X238").Value = -1
$ws.Range("Z238").Value = 1.3
$ws.Range("AA238").Value = -1
$ws.Range("AB238").Value = 1
$ws.Range("AC238").Value = 0.875
$ws.Range("AD238").Value = -1

# Row 239
$ws.Range("B239").Value = 6836277
$ws.Range("E239").Value = "CFR Cluj"
$ws.Range("F239").Value = "AFC Hermannstadt"
$ws.Range("G239").Value = 1
$ws.Range("K239").Value = "H"
$ws.Range("L239").Value = 1.7
$ws.Range("N239").Value = 5
$ws.Range("O239").Value = 1.65
$ws.Range("P239").Value = 3.5
$ws.Range("Q239").Value = 5.25
$ws.Range("R239").Value = -0.75
$ws.Range("S239").Value = 1.85
$ws.Range("T239").Value = 2
$ws.Range("U239").Value = 2.25
$ws.Range("V239").Value = 1.875
$ws.Range("W239").Value = 1.975
$ws.Range("X239").Value = 0.6499999999999999
$ws.Range("Y239").Value = -1
$ws.Range("AA239").Value = 0.425
$ws.Range("AD239").Value = 0.9750000000000001

# Row 309
$ws.Range("B309").Value = 8191462
$ws.Range("E309").Value = "CSM Politehnica Iasi"
$ws.Range("F309").Value = "Petrolul Ploiesti"
$ws.Range("G309").Value = 2
$ws.Range("H309").Value = 0
$ws.Range("J309").Value = 0
$ws.Range("K309").Value = "H"
$ws.Range("L309").Value = 2.1
$ws.Range("N309").Value = 3.1
$ws.Range("O309").Value = 1.8
$ws.Range("P309").Value = 3.2
$ws.Range("Q309").Value = 4.2
$ws.Range("R309").Value = -0.5
$ws.Range("S309").Value = 1.85
$ws.Range("T309").Value = 2
$ws.Range("V309").Value = 2.025
$ws.Range("W309").Value = 1.825
$ws.Range("X309").Value = 0.8
$ws.Range("Z309").Value = -1
$ws.Range("AA309").Value = 0.8500000000000001
$ws.Range("AB309").Value = -1
$ws.Range("AC309").Value = -0.5
$ws.Range("AD309").Value = 0.4125

# Row 310
$ws.Range("B310").Value = 8191463
$ws.Range("E310").Value = "Dinamo Bucharest"
$ws.Range("F310").Value = "ACS UTA Batrana Doamna"
$ws.Range("L310").Value = 1.833
$ws.Range("M310").Value = 3.4
$ws.Range("N310").Value = 3.6
$ws.Range("O310").Value = 1.5
$ws.Range("P310").Value = 4.333
$ws.Range("Q310").Value = 5
$ws.Range("R310").Value = -1
$ws.Range("S310").Value = 1.875
$ws.Range("T310").Value = 1.975
$ws.Range("U310").Value = 3
$ws.Range("V310").Value = 2.025
$ws.Range("W310").Value = 1.825
$ws.Range("X310").Value = 0.5
$ws.Range("AA310").Value = 0.875
$ws.Range("AC310").Value = -1
$ws.Range("AD310").Value = 0.825

# Row 311
$ws.Range("B311").Value = 8191523
$ws.Range("E311").Value = "Otelul Galati"
$ws.Range("F311").Value = "FC Botosani"
$ws.Range("I311").Value = 2
$ws.Range("L311").Value = 1.666
$ws.Range("M311").Value = 3.6
$ws.Range("N311").Value = 4.6
$ws.Range("O311").Value = 2.9
$ws.Range("P311").Value = 3.5
$ws.Range("Q311").Value = 2.2
$ws.Range("R311").Value = 0.25
$ws.Range("V311").Value = 1.875
$ws.Range("W311").Value = 1.975
$ws.Range("X311").Value = 1.9
$ws.Range("AD311").Value = 0.4875

# Row 312
$ws.Range("B312").Value = 8191475
$ws.Range("E312").Value = "FC U Craiova 1948"
$ws.Range("F312").Value = "AFC Hermannstadt"
$ws.Range("G312").Value = 1
$ws.Range("H312").Value = 3
$ws.Range("I312").Value = 0
$ws.Range("K312").Value = "A"
$ws.Range("L312").Value = 2.625
$ws.Range("M312").Value = 3.3
$ws.Range("N312").Value = 2.45
$ws.Range("O312").Value = 2.05
$ws.Range("P312").Value = 3.5
$ws.Range("Q312").Value = 3
$ws.Range("R312").Value = -0.25
$ws.Range("S312").Value = 1.85
$ws.Range("T312").Value = 2
$ws.Range("U312").Value = 2.25
$ws.Range("V312").Value = 1.825
$ws.Range("W312").Value = 2.025
$ws.Range("X312").Value = -1
$ws.Range("Z312").Value = 2
$ws.Range("AA312").Value = -1
$ws.Range("AB312").Value = 1
$ws.Range("AC312").Value = 0.825
$ws.Range("AD312").Value = -1

# Row 313
$ws.Range("B313").Value = 8191476
$ws.Range("E313").Value = "FC Voluntari"
$ws.Range("F313").Value = "Universitatea Cluj"
$ws.Range("G313").Value = 0
$ws.Range("H313").Value = 1
$ws.Range("J313").Value = 1
$ws.Range("L313").Value = 3.05
$ws.Range("N313").Value = 2.15
$ws.Range("O313").Value = 2.6
$ws.Range("P313").Value = 3.4
$ws.Range("Q313").Value = 2.4
$ws.Range("R313").Value = 0
$ws.Range("S313").Value = 2
$ws.Range("T313").Value = 1.85
$ws.Range("V313").Value = 2
$ws.Range("W313").Value = 1.85
$ws.Range("Z313").Value = 1.4
$ws.Range("AB313").Value = 0.8500000000000001
$ws.Range("AC313").Value = -1
$ws.Range("AD313").Value = 0.8500000000000001
